$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost/Nord coordinate values in row 3 to whole numbers
$ws.Range("Q3").Value = 798208
$ws.Range("R3").Value = 7232634

# Remove the Starttid (Z3) and Sluttid (AB3) values for row 3
$ws.Range("Z3").Value = $null
$ws.Range("AB3").Value = $null
